$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Un-merge the ranges that previously held one shared value per block.
#    (E2:E6, F2:F7, F8:F11) so every row can carry its own explicit value.
#    Un-merging alone preserves each cell's existing (center/center, etc.)
#    formatting, so no extra alignment work is required afterwards.
# ---------------------------------------------------------------------------
$ws.Range("E2:E6").UnMerge()
$ws.Range("F2:F7").UnMerge()
$ws.Range("F8:F11").UnMerge()

# ---------------------------------------------------------------------------
# 2. New header cells for the two columns that used to be merged, formatted
#    like the existing "Optimizations" header (D1).
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "FP"
$ws.Range("F1").Value = "LM"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:F1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-populate column E ("FP-16" / "FP-32") for every data row.
# ---------------------------------------------------------------------------
$ws.Range("E2:E6").Value = "FP-16"
$ws.Range("E12").Value = "FP-32"
$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("E12").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Re-populate column F ("roberta" / "distilbert", now lower-case) for
#    every data row.
# ---------------------------------------------------------------------------
$ws.Range("F2:F7").Value = "roberta"
$ws.Range("F8:F12").Value = "distilbert"
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("F12").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Add the new row 12 (wdc_cameras_small / cls_sep / 0.7946 / ...).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "wdc_cameras_small"
$ws.Range("B12").Value = "cls_sep"
$ws.Range("C12").Value = 0.7946
$ws.Range("D12").Value = "da: entry_swap - dk: None - summarize: False"

# ---------------------------------------------------------------------------
# 6. Selection bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("F8:F12").Select()
